$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B18:I18").Copy() | Out-Null
$ws.Range("B26:I26").PasteSpecial(-4122) | Out-Null
